# Update "想去人数" (column F, "number of people interested") counts.
# The same underlying events appear on both the "展览" sheet and the
# "全部类型" aggregate sheet, but the aggregate sheet has one extra row
# inserted above (a 演出/performance entry), so the row numbers for the
# later events are shifted down by one there.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 103
$ws1.Range("F7").Value  = 11658
$ws1.Range("F8").Value  = 4380
$ws1.Range("F17").Value = 5079
$ws1.Range("F18").Value = 62

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 103
$ws4.Range("F7").Value  = 11658
$ws4.Range("F8").Value  = 4380
$ws4.Range("F18").Value = 5079
$ws4.Range("F19").Value = 62

$wb.Save()
